$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.954.22"
$ws.Range("E2").Value = "  +1.17%  "

$ws.Range("D3").Value = "3.705.85"
$ws.Range("E3").Value = "  +3.99%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'245.04"
$ws.Range("E5").Value = "  +2.66%  "

$ws.Range("D6").Value = "'1.91"
$ws.Range("E6").Value = "  +18.88%  "

$ws.Range("D7").Value = "'671.39"
$ws.Range("E7").Value = "  +2.68%  "

$ws.Range("D8").Value = "'0.432"
$ws.Range("E8").Value = "  +7.21%  "

$ws.Range("D9").Value = "'1.13"
$ws.Range("E9").Value = "  +7.83%  "

$ws.Range("E10").Value = "  -0.04%  "

$ws.Range("D11").Value = "3.704.49"
$ws.Range("E11").Value = "  +4.03%  "

$ws.Range("D12").Value = "'45.43"
$ws.Range("E12").Value = "  +5.04%  "

$ws.Range("E13").Value = "  +2.19%  "

$ws.Range("D14").Value = "'6.63"
$ws.Range("E14").Value = "  +3.71%  "

$ws.Range("D15").Value = "4.396.38"
$ws.Range("E15").Value = "  +4.04%  "

$ws.Range("D16").Value = "'0.0000269"
$ws.Range("E16").Value = "  +4.94%  "

$ws.Range("D17").Value = "96.700.29"
$ws.Range("E17").Value = "  +0.98%  "

$ws.Range("D18").Value = "'9.01"
$ws.Range("E18").Value = "  +16.51%  "

$ws.Range("D19").Value = "3.706.16"
$ws.Range("E19").Value = "  +4.14%  "

$ws.Range("D20").Value = "'12.98"
$ws.Range("E20").Value = "  +3.02%  "

$ws.Range("D21").Value = "'18.62"
$ws.Range("E21").Value = "  +5.40%  "

$ws.Range("D22").Value = "'0.541"
$ws.Range("E22").Value = "  +5.65%  "

$ws.Range("D23").Value = "'518.07"
$ws.Range("E23").Value = "  +3.33%  "

$ws.Range("D24").Value = "'3.48"
$ws.Range("E24").Value = "  +3.24%  "

$ws.Range("D25").Value = "'0.0000210"
$ws.Range("E25").Value = "  +6.91%  "

$ws.Range("D26").Value = "'6.97"
$ws.Range("E26").Value = "  +0.53%  "

$ws.Range("D27").Value = "'102.11"
$ws.Range("E27").Value = "  +6.87%  "

$ws.Range("D28").Value = "'13.18"
$ws.Range("E28").Value = "  +3.84%  "

$ws.Range("D29").Value = "'0.168"
$ws.Range("E29").Value = "  +10.07%  "

$ws.Range("D30").Value = "'3.11"
$ws.Range("E30").Value = "  +3.95%  "

$ws.Range("D31").Value = "'12.15"
$ws.Range("E31").Value = "  +7.47%  "

$ws.Range("D32").Value = "'0.998"
$ws.Range("E32").Value = "  -0.13%  "

$ws.Range("E33").Value = "  +2.84%  "

$ws.Range("D34").Value = "'33.22"
$ws.Range("E34").Value = "  +6.40%  "

$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("D36").Value = "'0.594"
$ws.Range("E36").Value = "  +5.65%  "

$ws.Range("E37").Value = "  +7.33%  "

$ws.Range("D38").Value = "'8.85"
$ws.Range("E38").Value = "  +0.87%  "

$ws.Range("D39").Value = "'616.17"
$ws.Range("E39").Value = "  +0.81%  "

$ws.Range("D40").Value = "'42.78"
$ws.Range("E40").Value = "  +26.33%  "

$ws.Range("E41").Value = "  +9.60%  "

$ws.Range("D42").Value = "'0.968"
$ws.Range("E42").Value = "  +7.50%  "

$ws.Range("E43").Value = "  +8.88%  "

$ws.Range("E44").Value = "  -0.04%  "

$ws.Range("D45").Value = "'6.21"
$ws.Range("E45").Value = "  +9.41%  "

$ws.Range("D46").Value = "'0.0457"
$ws.Range("E46").Value = "  +9.01%  "

$ws.Range("D47").Value = "'0.434"
$ws.Range("E47").Value = "  +27.37%  "

$ws.Range("E48").Value = "  +2.67%  "

$ws.Range("D49").Value = "'23.62"
$ws.Range("E49").Value = "  +0.45%  "

$ws.Range("D50").Value = "'8.64"
$ws.Range("E50").Value = "  +6.17%  "

$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "'3.29"
$ws.Range("E51").Value = "  +4.90%  "

